# 每日学习.xlsx -- "Add files via upload" commit
# Adds two new diary entries (rows 28-29) below the existing row 25-27
# block, mirroring its layout (merged A:C note cell + a time-range cell in
# column E), and marks the newest time-range entry in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-home the centred formatting that rows 25-27 used onto the shared
#     style slot already used by rows 20-24, freeing it up for the new
#     rows 28-29 block (mirrors what Excel does when you insert a new
#     block before/under an old one that had its own duplicate style).
$ws.Range("A25:C27").HorizontalAlignment = -4108
$ws.Range("E25:E26").HorizontalAlignment = -4108

# --- New row 28: note text in merged A28:C29, first new time-range note in E28
$ws.Range("A28").Value = "这两天，一方面在推进度，另一方面，对之前的内容进行重读，编写，网上找讲解，感觉之前的内容，没有读太明白，最近代码写起来有点吃力"
$ws.Range("A28:C29").Merge()
$ws.Range("A28:C29").HorizontalAlignment = -4108
$ws.Range("E28").Value = "把之前内容重新看了看，进度推到4.9初始化"

# --- New row 29: second new time-range note in E29, highlighted in red
$ws.Range("E29").Value = "然后，看完后发现，不是代码问题，是我自己理解问题，之前没有理解getchar函数等等"
$ws.Range("E29").Font.Color = 255

# --- Move the active selection to where the user would land next (F34)
$ws.Range("F34").Select()
